# Natmi following Dr Hou advice
# Rebuild the LR-pair table for Ly86-Cd180 so every (sending cluster x
# target cluster) combination is present, including the new "sCs" cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "M1"
$ws.Cells.Item(2,2).Value = "Ly86"
$ws.Cells.Item(2,3).Value = "Cd180"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 126.7341283333333
$ws.Cells.Item(2,8).Value = 380.202385
$ws.Cells.Item(2,9).Value = 0.3866927553064318
$ws.Cells.Item(2,10).Value = 0.3866927553064317
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 1.257910333333333
$ws.Cells.Item(2,14).Value = 3.773731
$ws.Cells.Item(2,15).Value = 0.005038171381153627
$ws.Cells.Item(2,16).Value = 0.005038171381153626
$ws.Cells.Item(2,17).Value = 159.4201696164928
$ws.Cells.Item(2,18).Value = 1434.781526548435
$ws.Cells.Item(2,19).Value = 0.001948224373084307
$ws.Cells.Item(2,20).Value = 0.001948224373084306

# Row 3
$ws.Cells.Item(3,1).Value = "M1"
$ws.Cells.Item(3,2).Value = "Ly86"
$ws.Cells.Item(3,3).Value = "Cd180"
$ws.Cells.Item(3,4).Value = "M1"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 126.7341283333333
$ws.Cells.Item(3,8).Value = 380.202385
$ws.Cells.Item(3,9).Value = 0.3866927553064318
$ws.Cells.Item(3,10).Value = 0.3866927553064317
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 90.22063800000001
$ws.Cells.Item(3,14).Value = 270.661914
$ws.Cells.Item(3,15).Value = 0.3613509042067556
$ws.Cells.Item(3,16).Value = 0.3613509042067556
$ws.Cells.Item(3,17).Value = 11434.03391460721
$ws.Cells.Item(3,18).Value = 102906.3052314649
$ws.Cells.Item(3,19).Value = 0.1397317767801808
$ws.Cells.Item(3,20).Value = 0.1397317767801808

# Row 4
$ws.Cells.Item(4,1).Value = "M1"
$ws.Cells.Item(4,2).Value = "Ly86"
$ws.Cells.Item(4,3).Value = "Cd180"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 126.7341283333333
$ws.Cells.Item(4,8).Value = 380.202385
$ws.Cells.Item(4,9).Value = 0.3866927553064318
$ws.Cells.Item(4,10).Value = 0.3866927553064317
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 109.26873
$ws.Cells.Item(4,14).Value = 327.80619
$ws.Cells.Item(4,15).Value = 0.437642154415089
$ws.Cells.Item(4,16).Value = 0.437642154415089
$ws.Cells.Item(4,17).Value = 13848.07725064035
$ws.Cells.Item(4,18).Value = 124632.6952557632
$ws.Cells.Item(4,19).Value = 0.1692330505290136
$ws.Cells.Item(4,20).Value = 0.1692330505290136

# Row 5
$ws.Cells.Item(5,1).Value = "M1"
$ws.Cells.Item(5,2).Value = "Ly86"
$ws.Cells.Item(5,3).Value = "Cd180"
$ws.Cells.Item(5,4).Value = "Neutro"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 126.7341283333333
$ws.Cells.Item(5,8).Value = 380.202385
$ws.Cells.Item(5,9).Value = 0.3866927553064318
$ws.Cells.Item(5,10).Value = 0.3866927553064317
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 48.73834733333334
$ws.Cells.Item(5,14).Value = 146.215042
$ws.Cells.Item(5,15).Value = 0.1952063992103771
$ws.Cells.Item(5,16).Value = 0.1952063992103771
$ws.Cells.Item(5,17).Value = 6176.811965697241
$ws.Cells.Item(5,18).Value = 55591.30769127517
$ws.Cells.Item(5,19).Value = 0.07548490036410799
$ws.Cells.Item(5,20).Value = 0.07548490036410797

# Row 6
$ws.Cells.Item(6,1).Value = "M1"
$ws.Cells.Item(6,2).Value = "Ly86"
$ws.Cells.Item(6,3).Value = "Cd180"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 126.7341283333333
$ws.Cells.Item(6,8).Value = 380.202385
$ws.Cells.Item(6,9).Value = 0.3866927553064318
$ws.Cells.Item(6,10).Value = 0.3866927553064317
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.1903456666666667
$ws.Cells.Item(6,14).Value = 0.571037
$ws.Cells.Item(6,15).Value = 0.0007623707866246491
$ws.Cells.Item(6,16).Value = 0.0007623707866246489
$ws.Cells.Item(6,17).Value = 24.12329214702722
$ws.Cells.Item(6,18).Value = 217.109629323245
$ws.Cells.Item(6,19).Value = 0.0002948032600450173
$ws.Cells.Item(6,20).Value = 0.0002948032600450173

# Row 7
$ws.Cells.Item(7,1).Value = "M2"
$ws.Cells.Item(7,2).Value = "Ly86"
$ws.Cells.Item(7,3).Value = "Cd180"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 186.7837523333334
$ws.Cells.Item(7,8).Value = 560.351257
$ws.Cells.Item(7,9).Value = 0.569916918087593
$ws.Cells.Item(7,10).Value = 0.5699169180875929
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.257910333333333
$ws.Cells.Item(7,14).Value = 3.773731
$ws.Cells.Item(7,15).Value = 0.005038171381153627
$ws.Cells.Item(7,16).Value = 0.005038171381153626
$ws.Cells.Item(7,17).Value = 234.9572121588742
$ws.Cells.Item(7,18).Value = 2114.614909429867
$ws.Cells.Item(7,19).Value = 0.002871339106344187
$ws.Cells.Item(7,20).Value = 0.002871339106344186

# Row 8
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Ly86"
$ws.Cells.Item(8,3).Value = "Cd180"
$ws.Cells.Item(8,4).Value = "M1"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 186.7837523333334
$ws.Cells.Item(8,8).Value = 560.351257
$ws.Cells.Item(8,9).Value = 0.569916918087593
$ws.Cells.Item(8,10).Value = 0.5699169180875929
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 90.22063800000001
$ws.Cells.Item(8,14).Value = 270.661914
$ws.Cells.Item(8,15).Value = 0.3613509042067556
$ws.Cells.Item(8,16).Value = 0.3613509042067556
$ws.Cells.Item(8,17).Value = 16851.74930354733
$ws.Cells.Item(8,18).Value = 151665.7437319259
$ws.Cells.Item(8,19).Value = 0.2059399936736792
$ws.Cells.Item(8,20).Value = 0.2059399936736792

# Row 9
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Ly86"
$ws.Cells.Item(9,3).Value = "Cd180"
$ws.Cells.Item(9,4).Value = "M2"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 186.7837523333334
$ws.Cells.Item(9,8).Value = 560.351257
$ws.Cells.Item(9,9).Value = 0.569916918087593
$ws.Cells.Item(9,10).Value = 0.5699169180875929
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 109.26873
$ws.Cells.Item(9,14).Value = 327.80619
$ws.Cells.Item(9,15).Value = 0.437642154415089
$ws.Cells.Item(9,16).Value = 0.437642154415089
$ws.Cells.Item(9,17).Value = 20409.62340209787
$ws.Cells.Item(9,18).Value = 183686.6106188809
$ws.Cells.Item(9,19).Value = 0.249419667869462
$ws.Cells.Item(9,20).Value = 0.249419667869462

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Ly86"
$ws.Cells.Item(10,3).Value = "Cd180"
$ws.Cells.Item(10,4).Value = "Neutro"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 186.7837523333334
$ws.Cells.Item(10,8).Value = 560.351257
$ws.Cells.Item(10,9).Value = 0.569916918087593
$ws.Cells.Item(10,10).Value = 0.5699169180875929
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 48.73834733333334
$ws.Cells.Item(10,14).Value = 146.215042
$ws.Cells.Item(10,15).Value = 0.1952063992103771
$ws.Cells.Item(10,16).Value = 0.1952063992103771
$ws.Cells.Item(10,17).Value = 9103.531397445313
$ws.Cells.Item(10,18).Value = 81931.78257700781
$ws.Cells.Item(10,19).Value = 0.1112514294289545
$ws.Cells.Item(10,20).Value = 0.1112514294289544

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Ly86"
$ws.Cells.Item(11,3).Value = "Cd180"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 186.7837523333334
$ws.Cells.Item(11,8).Value = 560.351257
$ws.Cells.Item(11,9).Value = 0.569916918087593
$ws.Cells.Item(11,10).Value = 0.5699169180875929
$ws.Cells.Item(11,11).Value = 1
$ws.Cells.Item(11,12).Value = 0.3333333333333333
$ws.Cells.Item(11,13).Value = 0.1903456666666667
$ws.Cells.Item(11,14).Value = 0.571037
$ws.Cells.Item(11,15).Value = 0.0007623707866246491
$ws.Cells.Item(11,16).Value = 0.0007623707866246489
$ws.Cells.Item(11,17).Value = 35.55347786038989
$ws.Cells.Item(11,18).Value = 319.9813007435091
$ws.Cells.Item(11,19).Value = 0.000434488009153134
$ws.Cells.Item(11,20).Value = 0.0004344880091531338

# Row 12
$ws.Cells.Item(12,1).Value = "Neutro"
$ws.Cells.Item(12,2).Value = "Ly86"
$ws.Cells.Item(12,3).Value = "Cd180"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 14.22068333333333
$ws.Cells.Item(12,8).Value = 42.66205
$ws.Cells.Item(12,9).Value = 0.04339032660597528
$ws.Cells.Item(12,10).Value = 0.04339032660597528
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 1.257910333333333
$ws.Cells.Item(12,14).Value = 3.773731
$ws.Cells.Item(12,15).Value = 0.005038171381153627
$ws.Cells.Item(12,16).Value = 0.005038171381153626
$ws.Cells.Item(12,17).Value = 17.88834451206111
$ws.Cells.Item(12,18).Value = 160.99510060855
$ws.Cells.Item(12,19).Value = 0.0002186079017251335
$ws.Cells.Item(12,20).Value = 0.0002186079017251334

# Row 13
$ws.Cells.Item(13,1).Value = "Neutro"
$ws.Cells.Item(13,2).Value = "Ly86"
$ws.Cells.Item(13,3).Value = "Cd180"
$ws.Cells.Item(13,4).Value = "M1"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 14.22068333333333
$ws.Cells.Item(13,8).Value = 42.66205
$ws.Cells.Item(13,9).Value = 0.04339032660597528
$ws.Cells.Item(13,10).Value = 0.04339032660597528
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 90.22063800000001
$ws.Cells.Item(13,14).Value = 270.661914
$ws.Cells.Item(13,15).Value = 0.3613509042067556
$ws.Cells.Item(13,16).Value = 0.3613509042067556
$ws.Cells.Item(13,17).Value = 1282.9991231293
$ws.Cells.Item(13,18).Value = 11546.9921081637
$ws.Cells.Item(13,19).Value = 0.01567913375289562
$ws.Cells.Item(13,20).Value = 0.01567913375289561

# Row 14
$ws.Cells.Item(14,1).Value = "Neutro"
$ws.Cells.Item(14,2).Value = "Ly86"
$ws.Cells.Item(14,3).Value = "Cd180"
$ws.Cells.Item(14,4).Value = "M2"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 14.22068333333333
$ws.Cells.Item(14,8).Value = 42.66205
$ws.Cells.Item(14,9).Value = 0.04339032660597528
$ws.Cells.Item(14,10).Value = 0.04339032660597528
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 109.26873
$ws.Cells.Item(14,14).Value = 327.80619
$ws.Cells.Item(14,15).Value = 0.437642154415089
$ws.Cells.Item(14,16).Value = 0.437642154415089
$ws.Cells.Item(14,17).Value = 1553.8760075655
$ws.Cells.Item(14,18).Value = 13984.8840680895
$ws.Cells.Item(14,19).Value = 0.01898943601661338
$ws.Cells.Item(14,20).Value = 0.01898943601661338

# Row 15
$ws.Cells.Item(15,1).Value = "Neutro"
$ws.Cells.Item(15,2).Value = "Ly86"
$ws.Cells.Item(15,3).Value = "Cd180"
$ws.Cells.Item(15,4).Value = "Neutro"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 14.22068333333333
$ws.Cells.Item(15,8).Value = 42.66205
$ws.Cells.Item(15,9).Value = 0.04339032660597528
$ws.Cells.Item(15,10).Value = 0.04339032660597528
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 48.73834733333334
$ws.Cells.Item(15,14).Value = 146.215042
$ws.Cells.Item(15,15).Value = 0.1952063992103771
$ws.Cells.Item(15,16).Value = 0.1952063992103771
$ws.Cells.Item(15,17).Value = 693.0926036173445
$ws.Cells.Item(15,18).Value = 6237.833432556101
$ws.Cells.Item(15,19).Value = 0.008470069417314659
$ws.Cells.Item(15,20).Value = 0.008470069417314655

# Row 16
$ws.Cells.Item(16,1).Value = "Neutro"
$ws.Cells.Item(16,2).Value = "Ly86"
$ws.Cells.Item(16,3).Value = "Cd180"
$ws.Cells.Item(16,4).Value = "sCs"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 14.22068333333333
$ws.Cells.Item(16,8).Value = 42.66205
$ws.Cells.Item(16,9).Value = 0.04339032660597528
$ws.Cells.Item(16,10).Value = 0.04339032660597528
$ws.Cells.Item(16,11).Value = 1
$ws.Cells.Item(16,12).Value = 0.3333333333333333
$ws.Cells.Item(16,13).Value = 0.1903456666666667
$ws.Cells.Item(16,14).Value = 0.571037
$ws.Cells.Item(16,15).Value = 0.0007623707866246491
$ws.Cells.Item(16,16).Value = 0.0007623707866246489
$ws.Cells.Item(16,17).Value = 2.706845449538889
$ws.Cells.Item(16,18).Value = 24.36160904585
$ws.Cells.Item(16,19).Value = 0.00003307951742649781
$ws.Cells.Item(16,20).Value = 0.00003307951742649781
